# [AFG] added final excel sheets for Afghanistan
#
# 1) Remove the stray empty "INNING_NUMBER" cells at B7 / B14 on "ODI Batting".
# 2) Add two new worksheets at the end of the workbook:
#       "ODI Batting Extra"  (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
#                              PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH)
#       "ODI Bowling Extra"  (MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL)

$wb = $excel.ActiveWorkbook

# --- 1. Clear the two stray blank cells on "ODI Batting" ------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B7").ClearContents()
$odiBatting.Range("B14").ClearContents()

# --- 2. Add "ODI Batting Extra" as the 4th sheet ---------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "ODI Batting Extra"

# --- 3. Add "ODI Bowling Extra" as the 5th (last) sheet --------------------
$ws5 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "ODI Bowling Extra"

# --- Headers (bold / bordered / centered, matching the other sheets) -------
$ws4.Range("A1").Value = "MATCH_CODE"
$ws4.Range("B1").Value = "BATTING_POSITION"
$ws4.Range("C1").Value = "NUM_4"
$ws4.Range("D1").Value = "NUM_6"
$ws4.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws4.Range("F1").Value = "MAN_OF_MATCH"

$ws5.Range("A1").Value = "MATCH_CODE"
$ws5.Range("B1").Value = "MAIDEN_OVERS"
$ws5.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Copy the existing header formatting (bold, border, center/top align) from
# the "Player Info" header row so the new header cells look the same as the
# headers used across the rest of the workbook.
$headerSource = $wb.Worksheets.Item("Player Info").Range("A1")
$headerSource.Copy()
$ws4.Range("A1:F1").PasteSpecial(-4122)
$headerSource.Copy()
$ws5.Range("A1:C1").PasteSpecial(-4122)

# --- "ODI Batting Extra" data rows -----------------------------------------
$ws4.Range("A2").Value = "'3651"
$ws4.Range("F2").Value = "NO"

$ws4.Range("A3").Value = "'3652"
$ws4.Range("B3").Value = 4
$ws4.Range("C3").Value = "'2"
$ws4.Range("D3").Value = "'0"
$ws4.Range("E3").Value = "'15.06%"
$ws4.Range("F3").Value = "NO"

$ws4.Range("A4").Value = "'3721"
$ws4.Range("F4").Value = "NO"

$ws4.Range("A5").Value = "'3722"
$ws4.Range("B5").Value = 4
$ws4.Range("C5").Value = "'0"
$ws4.Range("D5").Value = "'0"
$ws4.Range("E5").Value = "'4.44%"
$ws4.Range("F5").Value = "NO"

$ws4.Range("A6").Value = "'3785"
$ws4.Range("B6").Value = 4
$ws4.Range("C6").Value = "'2"
$ws4.Range("D6").Value = "'0"
$ws4.Range("E6").Value = "'15.32%"
$ws4.Range("F6").Value = "NO"

$ws4.Range("A7").Value = "'4046"
$ws4.Range("F7").Value = "NO"

$ws4.Range("A8").Value = "'4092"
$ws4.Range("B8").Value = 4
$ws4.Range("C8").Value = "'2"
$ws4.Range("D8").Value = "'1"
$ws4.Range("E8").Value = "'22.27%"
$ws4.Range("F8").Value = "NO"

$ws4.Range("A9").Value = "'4093"
$ws4.Range("B9").Value = 5
$ws4.Range("C9").Value = "'0"
$ws4.Range("D9").Value = "'0"
$ws4.Range("E9").Value = "'11.36%"
$ws4.Range("F9").Value = "NO"

$ws4.Range("A10").Value = "'4097"
$ws4.Range("B10").Value = 5
$ws4.Range("C10").Value = "'1"
$ws4.Range("D10").Value = "'0"
$ws4.Range("E10").Value = "'10.17%"
$ws4.Range("F10").Value = "NO"

$ws4.Range("A11").Value = "'4129"
$ws4.Range("B11").Value = 5
$ws4.Range("C11").Value = "'2"
$ws4.Range("D11").Value = "'0"
$ws4.Range("E11").Value = "'9.31%"
$ws4.Range("F11").Value = "NO"

$ws4.Range("A12").Value = "'4131"
$ws4.Range("B12").Value = 5
$ws4.Range("C12").Value = "'0"
$ws4.Range("D12").Value = "'0"
$ws4.Range("F12").Value = "NO"

$ws4.Range("A13").Value = "'4132"
$ws4.Range("F13").Value = "NO"

$ws4.Range("A14").Value = "'4134"
$ws4.Range("B14").Value = 6
$ws4.Range("C14").Value = "'1"
$ws4.Range("D14").Value = "'0"
$ws4.Range("E14").Value = "'2.90%"
$ws4.Range("F14").Value = "NO"

$ws4.Range("A15").Value = "'4136"
$ws4.Range("B15").Value = 4
$ws4.Range("C15").Value = "'2"
$ws4.Range("D15").Value = "'0"
$ws4.Range("E15").Value = "'5.88%"
$ws4.Range("F15").Value = "NO"

$ws4.Range("A16").Value = "'4140"
$ws4.Range("B16").Value = 4
$ws4.Range("C16").Value = "'0"
$ws4.Range("D16").Value = "'0"
$ws4.Range("E16").Value = "'4.12%"
$ws4.Range("F16").Value = "NO"

$ws4.Range("A17").Value = "'4145"

# --- "ODI Bowling Extra" data rows -----------------------------------------
$ws5.Range("A2").Value = "'3722"
$ws5.Range("B2").Value = "'0"
$ws5.Range("C2").Value = "'"
